# ProductBacklog and recorded videos (Marco)
#
# Updates the "Product Backlog" sheet: renames the "Admin" role/pages to
# "Super Administrator" and fills in the three previously-empty backlog
# rows (Create User / View User / Edit User) for the super administrator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Row 18 - "Main Page": update the user-story text to mention "super admin"
# instead of "admin" (Notes column I is unchanged).
$ws.Range("H18").Value = "Given that I am a user whoever is doctor, patient or super admin, when I want to learn more about Neighborhood Doctors and create an account or log in my account, then I am able to log in to the system or choose whether I should create an account. "

# Row 21 - Sign Up page for the admin role is renamed to Super Administrator.
$ws.Range("C21").Value = "Super Administrator Sign Up Page"
$ws.Range("D21").Value = "Super Administrator"

# Row 24 - Log In page for the admin role is renamed to Super Administrator.
$ws.Range("C24").Value = "Super Administrator Log In Page"
$ws.Range("D24").Value = "Super Administrator"

# Row 25 - previously blank, now the "Create User" backlog item.
$ws.Range("C25").Value = "Super Administrator Create User Page"
$ws.Range("D25").Value = "Super Administrator"
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = "To Do "
$ws.Range("H25").Value = "Given that I am an admin, when I am going to create a new user, then I am able to add a new user to the system."
$ws.Range("I25").Value = "Access and add new information to the database."

# Row 26 - previously blank, now the "View User" backlog item.
$ws.Range("C26").Value = "Super Administrator View Page"
$ws.Range("D26").Value = "Super Administrator"
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = "To Do "
$ws.Range("H26").Value = "Given that I am an admin, when I am going to view user details, then I am able to view any user in the system using their unique ID."
$ws.Range("I26").Value = "Access user information from the database."

# Row 27 - previously blank, now the "Edit User" backlog item.
$ws.Range("C27").Value = "Super Administrator Edit User Page"
$ws.Range("D27").Value = "Super Administrator"
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = "To Do "
$ws.Range("H27").Value = "Given that I am an admin, when I am going to edit user details, then I am able to edit any user details and store it in the system."
$ws.Range("I27").Value = "Access and edit user information from the database."

# Leave the selection on the last cell that was touched, matching the
# author's final cursor position after entering the new rows.
$ws.Activate() | Out-Null
$ws.Range("B27").Select() | Out-Null

